$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build style s=5 (fillId=34, borderId=11) by copying formats from an
# existing "fillId=34 / borderId=10" cell (e.g. F2) and then removing the
# top/bottom border so only left/right remain.
$ws.Range("F2").Copy()
$ws.Range("G2").PasteSpecial(-4122)
$ws.Range("G2").Borders.Item(8).LineStyle = -4142
$ws.Range("G2").Borders.Item(9).LineStyle = -4142
$ws.Range("G2").Value = "ok"

$ws.Range("F3").Copy()
$ws.Range("G3").PasteSpecial(-4122)
$ws.Range("G3").Borders.Item(8).LineStyle = -4142
$ws.Range("G3").Borders.Item(9).LineStyle = -4142
$ws.Range("G3").Value = "ok"

$ws.Range("F5").Copy()
$ws.Range("G5").PasteSpecial(-4122)
$ws.Range("G5").Borders.Item(8).LineStyle = -4142
$ws.Range("G5").Borders.Item(9).LineStyle = -4142
$ws.Range("G5").Value = "ok"

$ws.Range("F6").Copy()
$ws.Range("G6").PasteSpecial(-4122)
$ws.Range("G6").Borders.Item(8).LineStyle = -4142
$ws.Range("G6").Borders.Item(9).LineStyle = -4142
$ws.Range("G6").Value = "ok"

$ws.Range("F7").Copy()
$ws.Range("G7").PasteSpecial(-4122)
$ws.Range("G7").Borders.Item(8).LineStyle = -4142
$ws.Range("G7").Borders.Item(9).LineStyle = -4142
$ws.Range("G7").Value = "ok"

$ws.Range("F9").Copy()
$ws.Range("G9").PasteSpecial(-4122)
$ws.Range("G9").Borders.Item(8).LineStyle = -4142
$ws.Range("G9").Borders.Item(9).LineStyle = -4142
$ws.Range("G9").Value = "ok"

$ws.Range("F10").Copy()
$ws.Range("G10").PasteSpecial(-4122)
$ws.Range("G10").Borders.Item(8).LineStyle = -4142
$ws.Range("G10").Borders.Item(9).LineStyle = -4142
$ws.Range("G10").Value = "ok des 74HC14"

$ws.Range("E8").Select()

Write-Host "done"
